# quarterly.xlsx update: roll the quarterly window forward one quarter
# (drop "1399/06", add "1401/12") and refresh the underlying figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Quarter headers (row 8 and row 24) -----------------------------------
$quarters = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)

$col = 5  # column E
foreach ($q in $quarters) {
    $ws.Cells.Item(8, $col).Value = $q
    $ws.Cells.Item(24, $col).Value = $q
    $col = $col + 1
}

# --- Data rows --------------------------------------------------------------
# Each row's ten quarterly figures (E:N) roll forward by one quarter: the
# oldest quarter's figure is dropped and a new figure is appended at N.

$rowData = @{
    14 = @(4741, 6221, 4315, 4747, 5826, 6623, 6824, 7138, 7670, 11563)
    15 = @(534, 392, 385, -27, 1888, 421, 482, 5932, 4936, 4097)
    16 = @(9468, 19493, 12155, 15864, 21703, 22898, 21794, 22877, 27324, 30947)
    17 = @(-8983, 107234, 125052, 155310, 106147, 196794, 193493, 231830, 208505, 281997)
    19 = @(201817, 2521346, 156930, 288539, 148530, 6299472, 1720634, 1865868, 2028317, 3268010)
    20 = @(207577, 2654686, 298837, 464433, 284094, 6526208, 1943227, 2133645, 2276752, 3596614)
    26 = @(4802, 158, 170, 173, 177, 5145, 5180, 5348, 5374, 5704)
    27 = @(156, 4903, 1832, 5027, 5190, 229, 224, 237, 243, 235)
}

foreach ($r in $rowData.Keys) {
    $values = $rowData[$r]
    $col = 5  # column E
    foreach ($v in $values) {
        $ws.Cells.Item($r, $col).Value = $v
        $col = $col + 1
    }
}
